$d = $word.ActiveDocument

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function New-PkgXml($bodyInner) {
    return '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document ' + $wNs + '><w:body>' + $bodyInner + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
}

# ---------------------------------------------------------------------------
# 1) Hardpoints paragraph
# ---------------------------------------------------------------------------
$hpPara = $d.Paragraphs(19)
$hpFull = $hpPara.Range
$hpTarget = $d.Range($hpFull.Start, $hpFull.End - 1)

$hpInner = '<w:p>' + `
  '<w:r><w:t>Many vehicles also have hardpoints, which represent things like the driver' + [char]0x2019 + 's seat/cabin, weapons, and other points of interest. These hardpoints can be targeted specifically instead of the vehicle as a whole, and as such have their own Armor Class and hit point pools. Damage to a hard point does not apply to the vehicl</w:t></w:r>' + `
  '<w:r><w:t>e as a whole</w:t></w:r>' + `
  '<w:r><w:t xml:space="preserve">, but can inflict </w:t></w:r>' + `
  '<w:proofErr w:type="spellStart"/>' + `
  '<w:r><w:t>affects</w:t></w:r>' + `
  '<w:proofErr w:type="spellEnd"/>' + `
  '<w:r><w:t xml:space="preserve"> on it such as reducing its speed, disabling a weapon, etc.</w:t></w:r>' + `
  '<w:r><w:t xml:space="preserve"> The specifics of a vehicle' + [char]0x2019 + 's hardpoints are explained in its stat block.</w:t></w:r>' + `
  '</w:p>'

$hpTarget.InsertXML((New-PkgXml $hpInner))

Write-Output "hardpoints done"

# ---------------------------------------------------------------------------
# 2) Creature capacity paragraph
# ---------------------------------------------------------------------------
$ccPara = $d.Paragraphs(31)
$ccFull = $ccPara.Range
$ccTarget = $d.Range($ccFull.Start, $ccFull.End - 1)

$ccInner = '<w:p><w:pPr><w:ind w:firstLine="720"/></w:pPr>' + `
  '<w:r><w:lastRenderedPageBreak/><w:t>Creature capacity describes how many</w:t></w:r>' + `
  '<w:r><w:t xml:space="preserve"> Medium</w:t></w:r>' + `
  '<w:r><w:t>-</w:t></w:r>' + `
  '<w:r><w:t>size</w:t></w:r>' + `
  '<w:r><w:t xml:space="preserve"> creatures can ride the vehicle comfortably</w:t></w:r>' + `
  '<w:r><w:t>; for Large creatures, halve the creature capacity.</w:t></w:r>' + `
  '<w:r><w:t xml:space="preserve"> More creatures can fit by squeezing or by clinging to the outside of the vehicle.</w:t></w:r>' + `
  '</w:p>'

$ccTarget.InsertXML((New-PkgXml $ccInner))

Write-Output "creature capacity done"

# ---------------------------------------------------------------------------
# 3) Cargo capacity paragraph
# ---------------------------------------------------------------------------
$cgPara = $d.Paragraphs(32)
$cgFull = $cgPara.Range
$cgTarget = $d.Range($cgFull.Start, $cgFull.End - 1)

$cgInner = '<w:p>' + `
  '<w:r><w:t>Cargo capacity specifies how much cargo the vehicle can carry.</w:t></w:r>' + `
  '<w:r><w:t xml:space="preserve"> For every creature under the </w:t></w:r>' + `
  '<w:r><w:t xml:space="preserve">creature </w:t></w:r>' + `
  '<w:r><w:t xml:space="preserve">capacity limit of a vehicle, that vehicle has room for an additional 50 lbs. of cargo capacity. For example, a </w:t></w:r>' + `
  '<w:proofErr w:type="spellStart"/>' + `
  '<w:r><w:t>Corvega</w:t></w:r>' + `
  '<w:proofErr w:type="spellEnd"/>' + `
  '<w:r><w:t xml:space="preserve"> Blitz (creature capacity of 5) with only a driver would have an additional </w:t></w:r>' + `
  '<w:r><w:t xml:space="preserve">4 x 50 = </w:t></w:r>' + `
  '<w:r><w:t>200 lbs. of cargo capacity available to it.</w:t></w:r>' + `
  '</w:p>'

$cgTarget.InsertXML((New-PkgXml $cgInner))

Write-Output "cargo capacity done"
